$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.19169999999999
$ws.Range("E4").Value = 13.58329999999999
$ws.Range("E7").Value = 12.0634
$ws.Range("E8").Value = 13.3915
$ws.Range("C11").Value = -13.51829999999999
$ws.Range("C12").Value = -14.62380000000002
$ws.Range("E12").Value = 11.6749
$ws.Range("E14").Value = 13.929
$ws.Range("C15").Value = -11.6487
$ws.Range("E22").Value = 12.5873
